$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing row, values replaced)
$ws.Range("B2").Value = "[1.66666666666667;1.08333333333333]"
$ws.Range("C2").Value = "'1.0"
$ws.Range("C2").Style = "Normal"

# Row 3 (new)
$ws.Range("A3").Value = "'2"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "[1.30555555555556;1.17361111111111]"
$ws.Range("C3").Value = "'0.276595744680851"
$ws.Range("C3").Style = "Normal"

# Row 4 (new)
$ws.Range("A4").Value = "'3"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "[1.27546296296296;1.18113425925926]"
$ws.Range("C4").Value = "'0.0235934664246824"
$ws.Range("C4").Style = "Normal"

# Row 5 (new)
$ws.Range("A5").Value = "'4"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "[1.27295524691358;1.18176118827161]"
$ws.Range("C5").Value = "'0.0019699954538566"
$ws.Range("C5").Style = "Normal"

# Row 6 (new)
$ws.Range("A6").Value = "'5"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "[1.27274627057613;1.18181343235597]"
$ws.Range("C6").Value = "'0.0001641932428166"
$ws.Range("C6").Style = "Normal"

# Row 7 (new)
$ws.Range("A7").Value = "'6"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "[1.27272885588134;1.18181778602966]"
$ws.Range("C7").Value = "'1.36829574554376e-05"
$ws.Range("C7").Style = "Normal"

# Row 8 (new)
$ws.Range("A8").Value = "'7"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "[1.27272740465678;1.18181814883581]"
$ws.Range("C8").Value = "'1.1402477547978e-06"
$ws.Range("C8").Style = "Normal"
